$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 79800
$ws.Range("J75").Value = 79800
$ws.Range("L75").Value = 79800
$ws.Range("N75").Value = -81672

$ws.Range("H78").Value = 79800
$ws.Range("J78").Value = 79800
$ws.Range("L78").Value = 239400
$ws.Range("N78").Value = -248760

$ws.Range("H97").Value = 300
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 300
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 900
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -1892

$ws.Range("H112").Value = 1588.7273
$ws.Range("I112").Value = 1165.3334
$ws.Range("J112").Value = 1747.5
$ws.Range("K112").Value = 3496.0002
$ws.Range("L112").Value = 5242.5
$ws.Range("M112").Value = -2388.0002
$ws.Range("N112").Value = -7458.5

$ws.Range("H138").Value = 3065.9546
$ws.Range("J138").Value = 2993.648
$ws.Range("L138").Value = 8980.944
$ws.Range("N138").Value = -19260.944

$ws.Range("H141").Value = 3815.9
$ws.Range("I141").Value = 3853.125
$ws.Range("K141").Value = 11559.375
$ws.Range("M141").Value = -6379.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1802.3889
$ws.Range("J2").Value = 3370.4285
$ws.Range("L2").Value = 3370.4285
$ws.Range("N2").Value = -3596.4285

$ws.Range("H32").Value = 5471.0312
$ws.Range("I32").Value = 5243.517
$ws.Range("J32").Value = 7670.3335
$ws.Range("K32").Value = 5243.517
$ws.Range("L32").Value = 7670.3335
$ws.Range("M32").Value = -4956.517
$ws.Range("N32").Value = -8244.333500000001

$ws.Range("H55").Value = 28999
$ws.Range("I55").Value = 29998
$ws.Range("K55").Value = 29998
$ws.Range("M55").Value = -29683

$ws.Range("H61").Value = 4760.923
$ws.Range("J61").Value = 7138.857
$ws.Range("L61").Value = 7138.857
$ws.Range("N61").Value = -7562.857

$ws.Range("H74").Value = 2985.7856
$ws.Range("I74").Value = 2998.889
$ws.Range("K74").Value = 2998.889
$ws.Range("M74").Value = -2124.889

$ws.Range("H77").Value = 2985.7856
$ws.Range("I77").Value = 2998.889
$ws.Range("K77").Value = 14994.445
$ws.Range("M77").Value = -10626.445

$ws.Range("H116").Value = 1802.3889
$ws.Range("J116").Value = 3370.4285
$ws.Range("L116").Value = 3370.4285
$ws.Range("N116").Value = -7958.4285

$ws.Range("H136").Value = 4760.923
$ws.Range("J136").Value = 7138.857
$ws.Range("L136").Value = 21416.571
$ws.Range("N136").Value = -26516.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1802.3889
$ws.Range("J3").Value = 3370.4285
$ws.Range("L3").Value = 3370.4285
$ws.Range("N3").Value = -3598.4285

$ws.Range("H75").Value = 25254.166
$ws.Range("I75").Value = 25254.166
$ws.Range("K75").Value = 25254.166
$ws.Range("M75").Value = -24318.166

$ws.Range("H78").Value = 25254.166
$ws.Range("I78").Value = 25254.166
$ws.Range("K78").Value = 75762.49800000001
$ws.Range("M78").Value = -71082.49800000001

$ws.Range("H107").Value = 2565.9048
$ws.Range("I107").Value = 2293.55
$ws.Range("K107").Value = 2293.55
$ws.Range("M107").Value = -373.5500000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 21198.54
$ws.Range("I22").Value = 2498.4285
$ws.Range("J22").Value = 43015.332
$ws.Range("K22").Value = 2498.4285
$ws.Range("L22").Value = 43015.332
$ws.Range("M22").Value = -2148.4285
$ws.Range("N22").Value = -43715.332

$ws.Range("H122").Value = 1411.5
$ws.Range("I122").Value = 1411.5
$ws.Range("K122").Value = 4234.5
$ws.Range("M122").Value = -1784.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 18333584
$ws.Range("I4").Value = 18333584
$ws.Range("K4").Value = 55000752
$ws.Range("M4").Value = -55000640

$ws.Range("H23").Value = 266.83334
$ws.Range("I23").Value = 124
$ws.Range("J23").Value = 409.66666
$ws.Range("K23").Value = 372
$ws.Range("L23").Value = 1228.99998
$ws.Range("M23").Value = -137
$ws.Range("N23").Value = -1698.99998

$ws.Range("H63").Value = 14459
$ws.Range("I63").Value = 14459
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 43377
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -42628
$ws.Range("N63").ClearContents()

$ws.Range("H64").Value = 14132.667
$ws.Range("J64").Value = 19950
$ws.Range("L64").Value = 59850
$ws.Range("N64").Value = -60390

$ws.Range("H66").Value = 14459
$ws.Range("I66").Value = 14459
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 130131
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -126387
$ws.Range("N66").ClearContents()

$ws.Range("H67").Value = 14132.667
$ws.Range("J67").Value = 19950
$ws.Range("L67").Value = 59850
$ws.Range("N67").Value = -61722

$ws.Range("H122").Value = 492.5
$ws.Range("J122").Value = 502
$ws.Range("L122").Value = 4518
$ws.Range("N122").Value = -9418

$ws.Range("H129").Value = 1812.7858
$ws.Range("I129").Value = 545.8889
$ws.Range("J129").Value = 4093.2
$ws.Range("K129").Value = 1637.6667
$ws.Range("L129").Value = 12279.6
$ws.Range("M129").Value = 3362.3333
$ws.Range("N129").Value = -22279.6

$ws.Range("H132").Value = 4251.5654
$ws.Range("I132").Value = 2769.2307
$ws.Range("J132").Value = 6178.6
$ws.Range("K132").Value = 24923.0763
$ws.Range("L132").Value = 55607.4
$ws.Range("M132").Value = -22393.0763
$ws.Range("N132").Value = -60667.4

$ws.Range("H134").Value = 14633.777
$ws.Range("I134").Value = 1500
$ws.Range("J134").Value = 18386.285
$ws.Range("K134").Value = 4500
$ws.Range("L134").Value = 55158.855
$ws.Range("M134").Value = 570
$ws.Range("N134").Value = -65298.855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 1663.3334
$ws.Range("J41").Value = 2000
$ws.Range("L41").Value = 2000
$ws.Range("N41").Value = -2710

$ws.Range("H102").Value = 1758.8
$ws.Range("I102").Value = 1527.2858
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 1527.2858
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = 94.71419999999989
$ws.Range("N102").Value = -8244

$ws.Range("H132").Value = 2982.2
$ws.Range("I132").Value = 2982.2
$ws.Range("K132").Value = 8946.599999999999
$ws.Range("M132").Value = -6416.599999999999

$ws.Range("H135").Value = 285657.5
$ws.Range("J135").Value = 285657.5
$ws.Range("L135").Value = 285657.5
$ws.Range("N135").Value = -295797.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1047.6471
$ws.Range("J22").Value = 1253.5
$ws.Range("L22").Value = 1253.5
$ws.Range("N22").Value = -1843.5

$ws.Range("H27").Value = 1047.6471
$ws.Range("J27").Value = 1253.5
$ws.Range("L27").Value = 1253.5
$ws.Range("N27").Value = -1467.5

$ws.Range("H40").Value = 3317.8572
$ws.Range("I40").Value = 2383.4443
$ws.Range("J40").Value = 4999.8
$ws.Range("K40").Value = 2383.4443
$ws.Range("L40").Value = 4999.8
$ws.Range("M40").Value = -2247.4443
$ws.Range("N40").Value = -5271.8

$ws.Range("H46").Value = 2596.8
$ws.Range("I46").Value = 1050.2
$ws.Range("K46").Value = 1050.2
$ws.Range("M46").Value = -862.2

$ws.Range("H55").Value = 306.70587
$ws.Range("I55").Value = 262.1
$ws.Range("J55").Value = 370.42856
$ws.Range("K55").Value = 262.1
$ws.Range("L55").Value = 370.42856
$ws.Range("M55").Value = -89.10000000000002
$ws.Range("N55").Value = -716.4285600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 2000000
$ws.Range("I8").Value = 2000000
$ws.Range("K8").Value = 2000000
$ws.Range("M8").Value = -1999860

$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

$ws.Range("H41").Value = 15162.5
$ws.Range("I41").Value = 11989.5
$ws.Range("K41").Value = 11989.5
$ws.Range("M41").Value = -11599.5

$ws.Range("H62").Value = 5799.7144
$ws.Range("I62").Value = 4333.3335
$ws.Range("K62").Value = 4333.3335
$ws.Range("M62").Value = -3709.3335

$ws.Range("H65").Value = 5799.7144
$ws.Range("I65").Value = 4333.3335
$ws.Range("K65").Value = 21666.6675
$ws.Range("M65").Value = -18546.6675

$ws.Range("H107").Value = 432.66666
$ws.Range("I107").Value = 99
$ws.Range("J107").Value = 599.5
$ws.Range("K107").Value = 297
$ws.Range("L107").Value = 1798.5
$ws.Range("M107").Value = 1623
$ws.Range("N107").Value = -5638.5

$ws.Range("H132").Value = 2652.4443
$ws.Range("I132").Value = 2484.625
$ws.Range("J132").Value = 3995
$ws.Range("K132").Value = 7453.875
$ws.Range("L132").Value = 11985
$ws.Range("M132").Value = -4923.875
$ws.Range("N132").Value = -17045

$ws.Range("H136").Value = 4277.8667
$ws.Range("I136").Value = 4442
$ws.Range("J136").Value = 3621.3333
$ws.Range("K136").Value = 13326
$ws.Range("L136").Value = 10863.9999
$ws.Range("M136").Value = -10776
$ws.Range("N136").Value = -15963.9999
